## Generate Report for Archive
##
## The status "Ready for handoff" has moved on to "In Translation" - update
## every cell that shows that status (the Overview roll-up columns for each
## locale, plus the per-locale "Status" column) and shrink the Status-style
## columns so they stay auto-fit to the new (shorter) text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update every occurrence of the status text -------------------------

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- Re-fit the columns that rendered the status text --------------------
# The text got shorter, so shrink the columns that display it back down to
# fit the new value (matches what AutoFit would do after the content edit).

$overview.Range("E:F").ColumnWidth = 12.5
$zhcn.Range("C:C").ColumnWidth = 12.5
$dede.Range("C:C").ColumnWidth = 12.5
